$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.892066666666667
$ws.Range("H2").Value = 23.6762
$ws.Range("I2").Value = 0.1739002798877711
$ws.Range("J2").Value = 0.1739002798877711
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.4394373333333333
$ws.Range("N2").Value = 1.318312
$ws.Range("O2").Value = 0.2944741752765458
$ws.Range("P2").Value = 0.2944741752765458
$ws.Range("Q2").Value = 3.468068730488889
$ws.Range("R2").Value = 31.2126185744
$ws.Range("S2").Value = 0.05120914150031188
$ws.Range("T2").Value = 0.05120914150031188

$ws.Range("G3").Value = 7.892066666666667
$ws.Range("H3").Value = 23.6762
$ws.Range("I3").Value = 0.1739002798877711
$ws.Range("J3").Value = 0.1739002798877711
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.6503416666666667
$ws.Range("N3").Value = 1.951025
$ws.Range("O3").Value = 0.4358046333636673
$ws.Range("P3").Value = 0.4358046333636673
$ws.Range("Q3").Value = 5.132539789444445
$ws.Range("R3").Value = 46.19285810500001
$ws.Range("S3").Value = 0.0757865477183292
$ws.Range("T3").Value = 0.0757865477183292

$ws.Range("G4").Value = 7.892066666666667
$ws.Range("H4").Value = 23.6762
$ws.Range("I4").Value = 0.1739002798877711
$ws.Range("J4").Value = 0.1739002798877711
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1607546666666667
$ws.Range("N4").Value = 0.482264
$ws.Range("O4").Value = 0.1077243426939663
$ws.Range("P4").Value = 0.1077243426939663
$ws.Range("Q4").Value = 1.268686546311111
$ws.Range("R4").Value = 11.4181789168
$ws.Range("S4").Value = 0.01873329334520691
$ws.Range("T4").Value = 0.01873329334520691

$ws.Range("G5").Value = 7.892066666666667
$ws.Range("H5").Value = 23.6762
$ws.Range("I5").Value = 0.1739002798877711
$ws.Range("J5").Value = 0.1739002798877711
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2417443333333333
$ws.Range("N5").Value = 0.725233
$ws.Range("O5").Value = 0.1619968486658205
$ws.Range("P5").Value = 0.1619968486658205
$ws.Range("Q5").Value = 1.907862394955556
$ws.Range("R5").Value = 17.1707615546
$ws.Range("S5").Value = 0.02817129732392308
$ws.Range("T5").Value = 0.02817129732392308

$ws.Range("G6").Value = 17.317702
$ws.Range("H6").Value = 51.95310600000001
$ws.Range("I6").Value = 0.3815924715300191
$ws.Range("J6").Value = 0.3815924715300191
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.4394373333333333
$ws.Range("N6").Value = 1.318312
$ws.Range("O6").Value = 0.2944741752765458
$ws.Range("P6").Value = 0.2944741752765458
$ws.Range("Q6").Value = 7.610044786341333
$ws.Range("R6").Value = 68.490403077072
$ws.Range("S6").Value = 0.1123691283455412
$ws.Range("T6").Value = 0.1123691283455412

$ws.Range("G7").Value = 17.317702
$ws.Range("H7").Value = 51.95310600000001
$ws.Range("I7").Value = 0.3815924715300191
$ws.Range("J7").Value = 0.3815924715300191
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.6503416666666667
$ws.Range("N7").Value = 1.951025
$ws.Range("O7").Value = 0.4358046333636673
$ws.Range("P7").Value = 0.4358046333636673
$ws.Range("Q7").Value = 11.26242318151667
$ws.Range("R7").Value = 101.36180863365
$ws.Range("S7").Value = 0.1662997671494756
$ws.Range("T7").Value = 0.1662997671494756

$ws.Range("G8").Value = 17.317702
$ws.Range("H8").Value = 51.95310600000001
$ws.Range("I8").Value = 0.3815924715300191
$ws.Range("J8").Value = 0.3815924715300191
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.1607546666666667
$ws.Range("N8").Value = 0.482264
$ws.Range("O8").Value = 0.1077243426939663
$ws.Range("P8").Value = 0.1077243426939663
$ws.Range("Q8").Value = 2.783901412442667
$ws.Range("R8").Value = 25.05511271198401
$ws.Range("S8").Value = 0.04110679817253737
$ws.Range("T8").Value = 0.04110679817253737

$ws.Range("G9").Value = 17.317702
$ws.Range("H9").Value = 51.95310600000001
$ws.Range("I9").Value = 0.3815924715300191
$ws.Range("J9").Value = 0.3815924715300191
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.2417443333333333
$ws.Range("N9").Value = 0.725233
$ws.Range("O9").Value = 0.1619968486658205
$ws.Range("P9").Value = 0.1619968486658205
$ws.Range("Q9").Value = 4.186456324855333
$ws.Range("R9").Value = 37.678106923698
$ws.Range("S9").Value = 0.06181677786246493
$ws.Range("T9").Value = 0.06181677786246494

$ws.Range("G10").Value = 7.716272666666666
$ws.Range("H10").Value = 23.148818
$ws.Range("I10").Value = 0.1700266904854272
$ws.Range("J10").Value = 0.1700266904854272
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.4394373333333333
$ws.Range("N10").Value = 1.318312
$ws.Range("O10").Value = 0.2944741752765458
$ws.Range("P10").Value = 0.2944741752765458
$ws.Range("Q10").Value = 3.390818283912888
$ws.Range("R10").Value = 30.51736455521599
$ws.Range("S10").Value = 0.0500684694556967
$ws.Range("T10").Value = 0.05006846945569671

$ws.Range("G11").Value = 7.716272666666666
$ws.Range("H11").Value = 23.148818
$ws.Range("I11").Value = 0.1700266904854272
$ws.Range("J11").Value = 0.1700266904854272
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.6503416666666667
$ws.Range("N11").Value = 1.951025
$ws.Range("O11").Value = 0.4358046333636673
$ws.Range("P11").Value = 0.4358046333636673
$ws.Range("Q11").Value = 5.018213626494444
$ws.Range("R11").Value = 45.16392263845
$ws.Range("S11").Value = 0.07409841950903935
$ws.Range("T11").Value = 0.07409841950903935

$ws.Range("G12").Value = 7.716272666666666
$ws.Range("H12").Value = 23.148818
$ws.Range("I12").Value = 0.1700266904854272
$ws.Range("J12").Value = 0.1700266904854272
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.1607546666666667
$ws.Range("N12").Value = 0.482264
$ws.Range("O12").Value = 0.1077243426939663
$ws.Range("P12").Value = 0.1077243426939663
$ws.Range("Q12").Value = 1.240426840439111
$ws.Range("R12").Value = 11.163841563952
$ws.Range("S12").Value = 0.0183160134729731
$ws.Range("T12").Value = 0.0183160134729731

$ws.Range("G13").Value = 7.716272666666666
$ws.Range("H13").Value = 23.148818
$ws.Range("I13").Value = 0.1700266904854272
$ws.Range("J13").Value = 0.1700266904854272
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.2417443333333333
$ws.Range("N13").Value = 0.725233
$ws.Range("O13").Value = 0.1619968486658205
$ws.Range("P13").Value = 0.1619968486658205
$ws.Range("Q13").Value = 1.865365191621555
$ws.Range("R13").Value = 16.788286724594
$ws.Range("S13").Value = 0.02754378804771806
$ws.Range("T13").Value = 0.02754378804771807

$ws.Range("G14").Value = 12.45667266666667
$ws.Range("H14").Value = 37.370018
$ws.Range("I14").Value = 0.2744805580967825
$ws.Range("J14").Value = 0.2744805580967826
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.4394373333333333
$ws.Range("N14").Value = 1.318312
$ws.Range("O14").Value = 0.2944741752765458
$ws.Range("P14").Value = 0.2944741752765458
$ws.Range("Q14").Value = 5.473927018846222
$ws.Range("R14").Value = 49.265343169616
$ws.Range("S14").Value = 0.08082743597499606
$ws.Range("T14").Value = 0.08082743597499607

$ws.Range("G15").Value = 12.45667266666667
$ws.Range("H15").Value = 37.370018
$ws.Range("I15").Value = 0.2744805580967825
$ws.Range("J15").Value = 0.2744805580967826
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.6503416666666667
$ws.Range("N15").Value = 1.951025
$ws.Range("O15").Value = 0.4358046333636673
$ws.Range("P15").Value = 0.4358046333636673
$ws.Range("Q15").Value = 8.101093263161113
$ws.Range("R15").Value = 72.90983936845001
$ws.Range("S15").Value = 0.1196198989868231
$ws.Range("T15").Value = 0.1196198989868231

$ws.Range("G16").Value = 12.45667266666667
$ws.Range("H16").Value = 37.370018
$ws.Range("I16").Value = 0.2744805580967825
$ws.Range("J16").Value = 0.2744805580967826
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.1607546666666667
$ws.Range("N16").Value = 0.482264
$ws.Range("O16").Value = 0.1077243426939663
$ws.Range("P16").Value = 0.1077243426939663
$ws.Range("Q16").Value = 2.002468262305778
$ws.Range("R16").Value = 18.022214360752
$ws.Range("S16").Value = 0.02956823770324894
$ws.Range("T16").Value = 0.02956823770324894

$ws.Range("G17").Value = 12.45667266666667
$ws.Range("H17").Value = 37.370018
$ws.Range("I17").Value = 0.2744805580967825
$ws.Range("J17").Value = 0.2744805580967826
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.2417443333333333
$ws.Range("N17").Value = 0.725233
$ws.Range("O17").Value = 0.1619968486658205
$ws.Range("P17").Value = 0.1619968486658205
$ws.Range("Q17").Value = 3.011330029354889
$ws.Range("R17").Value = 27.101970264194
$ws.Range("S17").Value = 0.04446498543171444
$ws.Range("T17").Value = 0.04446498543171445
